$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("WT")
$ws2 = $wb.Worksheets.Item("R183W")

# --- Sheet "R183W": insert a new blank row above row 12 ---
# This shifts the old rows 12-21 (the "ATUX (4,5 uM)" / "DMSO/DMSO" blocks) down to 13-22,
# matching the new dataset layout (dimension grows from K21 to K22).
$ws2.Rows.Item(12).Insert() | Out-Null

# --- Column K width/visibility on both sheets ---
# WT: column K was hidden (width 0) -> now shown with an explicit (best-fit-like) width.
$ws1.Columns.Item(11).Hidden = $false
$ws1.Columns.Item(11).ColumnWidth = 11.6

# R183W: gains its own explicit column width for column K as well.
$ws2.Columns.Item(11).ColumnWidth = 11.6

# --- Selections ---
$ws1.Range("A2").Select() | Out-Null
$ws2.Range("K10").Select() | Out-Null

# --- Active sheet / tab ---
# R183W becomes the active (selected) tab, replacing WT.
$ws2.Activate() | Out-Null
